$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.782044172286987
$ws.Range("B1").Value = 5.204721450805664
$ws.Range("C1").Value = 6.814132213592529
$ws.Range("D1").Value = 10.66283798217773
$ws.Range("E1").Value = 5.478264331817627
